$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.566.88"
$ws.Range("E2").Value = "  +2.85%  "

$ws.Range("D3").Value = "1.848.77"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("D4").Value = "1.034"
$ws.Range("E4").Value = "  +3.33%  "

$ws.Range("D5").Value = "320.78"
$ws.Range("E5").Value = "  +3.90%  "

$ws.Range("D6").Value = "1.029"
$ws.Range("E6").Value = "  +2.79%  "

$ws.Range("D7").Value = "0.4372"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("D8").Value = "0.3746"
$ws.Range("E8").Value = "  +1.70%  "

$ws.Range("D9").Value = "0.07389"
$ws.Range("E9").Value = "  +2.67%  "

$ws.Range("D10").Value = "0.8760"
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("D11").Value = "21.47"
$ws.Range("E11").Value = "  +3.05%  "

$ws.Range("D12").Value = "1.872.40"
$ws.Range("E12").Value = "  -3.44%  "

$ws.Range("D13").Value = "5.502"
$ws.Range("E13").Value = "  +3.07%  "

$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").Value = "0.07171"
$ws.Range("E15").Value = "  +4.15%  "

$ws.Range("D16").Value = "82.90"
$ws.Range("E16").Value = "  +3.12%  "

$ws.Range("D17").Value = "1.035"
$ws.Range("E17").Value = "  +3.36%  "

$ws.Range("D18").Value = "0.000009021"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("D19").Value = "1.028"
$ws.Range("E19").Value = "  +2.68%  "

$ws.Range("D20").Value = "15.42"
$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("D21").Value = "27.625.73"
$ws.Range("E21").Value = "  +2.96%  "

$ws.Range("D22").Value = "5.266"
$ws.Range("E22").Value = "  +1.16%  "

$ws.Range("D23").Value = "11.23"
$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("D24").Value = "2.066.25"
$ws.Range("E24").Value = "  -4.35%  "

$ws.Range("D25").Value = "157.53"
$ws.Range("E25").Value = "  +2.67%  "

$ws.Range("D26").Value = "1.934"
$ws.Range("E26").Value = "  +3.63%  "

$ws.Range("D27").Value = "18.75"
$ws.Range("E27").Value = "  +2.67%  "

$ws.Range("D28").Value = "5.282"
$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("D29").Value = "1.945"
$ws.Range("E29").Value = "  +2.41%  "

$ws.Range("D30").Value = "116.24"
$ws.Range("E30").Value = "  +0.90%  "

$ws.Range("D31").Value = "0.09072"
$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("D32").Value = "1.206"
$ws.Range("E32").Value = "  +3.53%  "

$ws.Range("D33").Value = "0.7669"
$ws.Range("E33").Value = "  +1.32%  "

$ws.Range("D34").Value = "4.515"
$ws.Range("E34").Value = "  +2.11%  "

$ws.Range("D35").Value = "2.878"
$ws.Range("E35").Value = "  +3.95%  "

$ws.Range("D36").Value = "1.030"
$ws.Range("E36").Value = "  +2.67%  "

$ws.Range("D37").Value = "1.151"
$ws.Range("E37").Value = "  +2.64%  "

$ws.Range("D38").Value = "0.01983"
$ws.Range("E38").Value = "  +3.14%  "

$ws.Range("E39").Value = "  +1.17%  "

$ws.Range("D40").Value = "0.5189"
$ws.Range("E40").Value = "  +2.28%  "

$ws.Range("D41").Value = "2.810"
$ws.Range("E41").Value = "  +6.12%  "

$ws.Range("D42").Value = "0.1672"
$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("D43").Value = "6.721"
$ws.Range("E43").Value = "  +3.06%  "

$ws.Range("D44").Value = "8.577"
$ws.Range("E44").Value = "  +3.66%  "

$ws.Range("D45").Value = "108.82"
$ws.Range("E45").Value = "  +2.38%  "

$ws.Range("E46").Value = "  +1.88%  "

$ws.Range("D47").Value = "1.721"
$ws.Range("E47").Value = "  +4.18%  "

$ws.Range("D48").Value = "0.4654"
$ws.Range("E48").Value = "  +2.53%  "

$ws.Range("D49").Value = "0.06389"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("D50").Value = "1.884"
$ws.Range("E50").Value = "  +4.38%  "

$ws.Range("D51").Value = "39.54"
$ws.Range("E51").Value = "  +5.86%  "
